$d = $word.ActiveDocument

# Locate the start of the "Prime fasi di test dell'interfaccia" section (the
# Titolo2 heading that begins the block being removed) using Find, which
# sidesteps any special-character (curly apostrophe) matching issues.
$findRange = $d.Content
$findRange.Find.Execute("Prime fasi di test dell", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$startPara = $findRange.Paragraphs.First

$lastPara = $d.Paragraphs.Last

# Delete everything from the start of the heading paragraph up to (but not
# including) the start of the very last paragraph in the document body.
# (Deleting a range that partially overlaps the final paragraph is
# unreliable, so the last paragraph is cleared out separately below.)
$r = $d.Range($startPara.Range.Start, $lastPara.Range.Start)
$r.Delete()

# The former last paragraph is now the sole remaining paragraph of the
# removed block. Clear its text (excluding its own paragraph mark)...
$remaining = $d.Paragraphs.Last
$rText = $d.Range($remaining.Range.Start, $remaining.Range.End - 1)
$rText.Delete()

# ...and strip its formatting/list numbering so it becomes a plain, unstyled
# empty paragraph.
$remaining2 = $d.Paragraphs.Last
$remaining2.Style = "Normal"
